$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Header row: add I1 ("I0") and J1 ("IF"), matching the style of H1 ---
$ws.Range("H1").Copy()
$ws.Range("I1:J1").PasteSpecial(-4122)
$ws.Range("I1").Value = "I0"
$ws.Range("J1").Value = "IF"

# --- Rows 2-15: I column is always 1, J column duplicates H column ---
$ws.Range("I2:I15").Value = 1
$ws.Range("H2:H15").Copy()
$ws.Range("J2:J15").PasteSpecial(-4163)

# --- Rows 16-25: explicit I/J values ---
$ws.Range("I16").Value = 7
$ws.Range("J16").Value = 8

$ws.Range("I17").Value = 6
$ws.Range("J17").Value = 9

$ws.Range("I18").Value = 5
$ws.Range("J18").Value = 7

$ws.Range("I19").Value = 6
$ws.Range("J19").Value = 6

$ws.Range("I20").Value = 6
$ws.Range("J20").Value = 8

$ws.Range("I21").Value = 6
$ws.Range("J21").Value = 7

$ws.Range("I22").Value = 6
$ws.Range("J22").Value = 7

$ws.Range("I23").Value = 4
$ws.Range("J23").Value = 6

$ws.Range("I24").Value = 4
$ws.Range("J24").Value = 5

$ws.Range("I25").Value = 2
$ws.Range("J25").Value = 2
